$d = $word.ActiveDocument

$hdr = $d.Sections.Item(1).Headers.Item(1)
$r = $hdr.Range.Duplicate
$r.Find.Execute("November")
Write-Output "start=$($r.Start) end=$($r.End) text=$($r.Text)"

# Try to set End to narrow the range to "Nov"
$rNov = $r.Duplicate
$rNov.End = $rNov.Start + 3
Write-Output "rNov text=$($rNov.Text)"
$rNov.Text = "Dec"

Write-Output $hdr.Range.Text
